# Append 11 new days of restaurant analytics data (rows 355-365) to Sheet1.
# Columns: A=day, B=idx, C=date, D=count, E=tavg, F=prcp, G=snow, H=wspd

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row data: day, idx, date(serial), count(optional/$null), tavg, prcp, snow, wspd
$rows = @(
    @{ r=355; A=2; B=354; C=45726; D=52;    E=55.2; F=0;    G=0;   H=11.1 },
    @{ r=356; A=3; B=355; C=45727; D=95;    E=45.1; F=0;    G=0;   H=12   },
    @{ r=357; A=4; B=356; C=45728; D=$null; E=43.9; F=0;    G=0;   H=8.1  },
    @{ r=358; A=5; B=357; C=45729; D=$null; E=49.6; F=0;    G=0;   H=7.5  },
    @{ r=359; A=6; B=358; C=45730; D=$null; E=58.5; F=0.15; G=0;   H=12.6 },
    @{ r=360; A=7; B=359; C=45731; D=$null; E=54.1; F=0.25; G=0;   H=18   },
    @{ r=361; A=1; B=360; C=45732; D=$null; E=40.5; F=0.05; G=1;   H=12.6 },
    @{ r=362; A=2; B=361; C=45733; D=$null; E=39.4; F=0;    G=0;   H=10.5 },
    @{ r=363; A=3; B=362; C=45734; D=$null; E=42.3; F=0;    G=0;   H=11.1 },
    @{ r=364; A=4; B=363; C=45735; D=$null; E=43.7; F=0;    G=0;   H=11   },
    @{ r=365; A=5; B=364; C=45736; D=$null; E=44.8; F=0.8;  G=0.1; H=11.2 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    if ($null -ne $row.D) {
        $ws.Cells.Item($r, 4).Value = $row.D
    }
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
}

# The date column (C) uses a date number-format style elsewhere in the
# column (s="3", numFmtId 14). New cells default to the "General" style,
# so copy the format from the last pre-existing date cell (C354) down
# across the newly written date cells, re-using that same style.
$null = $ws.Cells.Item(354, 3).Copy()
$null = $ws.Range("C355:C365").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the view: keep the header-row freeze, scroll the frozen pane so
# row 350 is the top visible data row, and select D363 as the active cell
# (matching the final editing position recorded in the workbook).
$null = $ws.Range("A350").Select()
$null = $ws.Range("D363").Select()
